$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Row 6 (Rusia) ---
$ws.Range("B6").Value = 654405
$ws.Range("C6").Value = 6556
$ws.Range("D6").Value = 422931
$ws.Range("E6").Value = 221938
$ws.Range("G6").Value = 216
$ws.Range("H6").Value = 9536

# --- Row 38 (Singapur) ---
$ws.Range("B38").Value = 44153
$ws.Range("C38").Value = 246
$ws.Range("E38").Value = 5627

# --- Row 46 (Afganistan) ---
$ws.Range("B46").Value = 31836
$ws.Range("C46").Value = 319
$ws.Range("D46").Value = 15651
$ws.Range("E46").Value = 15411

# --- Row 50 (now Armenia, new data, country moves above Nigeria) ---
$ws.Range("A50").Value = "Armenia"
$ws.Range("B50").Value = 26065
$ws.Range("C50").Value = 523
$ws.Range("D50").Value = 14563
$ws.Range("E50").Value = 11049
$ws.Range("G50").Value = 10
$ws.Range("H50").Value = 453

# --- Row 51 (now Nigeria, reusing the old Nigeria data) ---
$ws.Range("A51").Value = "Nigeria"
$ws.Range("B51").Value = 25694
$ws.Range("C51").Value = 0
$ws.Range("D51").Value = 9746
$ws.Range("E51").Value = 15358
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 590

# --- Row 69 (Chequia) ---
$ws.Range("B69").Value = 11960
$ws.Range("C69").Value = 6
$ws.Range("D69").Value = 7776
$ws.Range("E69").Value = 3835

# --- Row 83 (Tayikistan) ---
$ws.Range("B83").Value = 5954
$ws.Range("C83").Value = 54
$ws.Range("D83").Value = 4568
$ws.Range("E83").Value = 1334

# --- Row 94 (Hungria) ---
$ws.Range("B94").Value = 4157
$ws.Range("C94").Value = 2
$ws.Range("D94").Value = 2714
$ws.Range("E94").Value = 858

# --- Row 118 (Eslovaquia) ---
$ws.Range("B118").Value = 1687
$ws.Range("C118").Value = 20
$ws.Range("D118").Value = 1466
$ws.Range("E118").Value = 193

# --- Row 131 (Letonia) ---
$ws.Range("B131").Value = 1121
$ws.Range("C131").Value = 3
$ws.Range("E131").Value = 117

# --- Row 157 (Taiwan) ---
$ws.Range("D157").Value = 438
$ws.Range("E157").Value = 2

# --- Row 177 (Camboya) ---
$ws.Range("D177").Value = 131
$ws.Range("E177").Value = 10

# --- Row 193 (Islas Turcas y Caicos) ---
$ws.Range("B193").Value = 42
$ws.Range("C193").Value = 1
$ws.Range("E193").Value = 29

# --- Update "last updated" timestamp caption in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 1 de Julio de 2020 a las 10:04"
